# Insert a new data row at row 23 (shifts existing rows 23..118 down to 24..119,
# matching the weekly refresh described by the commit message "Fruta / hortaliza, semanal").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(23).Insert()

# Populate the newly inserted row with the new daily reading.
$ws.Cells.Item(23, 1).Value = 7
$ws.Cells.Item(23, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(23, 3).Value = "Ñuble"
$ws.Cells.Item(23, 4).Value = 45030
$ws.Cells.Item(23, 5).Value = 16
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100108
$ws.Cells.Item(23, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(23, 9).Value = 100108002
$ws.Cells.Item(23, 10).Value = "Mango"
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 70
$ws.Cells.Item(23, 14).Value = 8000
$ws.Cells.Item(23, 15).Value = 9000
$ws.Cells.Item(23, 16).Value = 8571
$ws.Cells.Item(23, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(23, 18).Value = "Perú"
$ws.Cells.Item(23, 19).Value = 2143
$ws.Cells.Item(23, 20).Value = 4
